# New crime data collected - weekly CompStat update (121st Precinct)
# Updates: report header (volume/number, week-of dates) and the Crime
# Complaints table (rows 15-27, 30) with refreshed weekly/28-day/YTD/2-year
# figures. A handful of cells flip between the literal text placeholders
# ("0" / "***.*", used when 2-Year % change is undefined) and real numbers
# as the underlying counts move away from/into that edge case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  41" -> "...42", and the week-of date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# ---------------------------------------------------------------------
# Helper cells used as "format donors" for PasteSpecial(xlPasteFormats):
#   Row 22 columns C:N are all styled as the shared "text placeholder"
#   style (s=14) - used whenever we need to turn a numeric cell into a
#   text "0" / "***.*" cell.
#   D16 carries the plain numeric style (s=15) - used when turning the
#   one remaining text cell (C16) back into a number.
# ---------------------------------------------------------------------

function Set-TextPlaceholder($ref, $text, $formatDonorRef) {
    # Write the literal text first (apostrophe forces text, not a number),
    # THEN paste the donor's number format over it - if done in the other
    # order, assigning .Value re-triggers Excel's "looks like text" auto
    # number-format (@) and clobbers the style we want to keep (s=14).
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($formatDonorRef).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumberFromText($ref, $value, $formatDonorRef) {
    $ws.Range($formatDonorRef).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $value
}

# ---------------------------------------------------------------------
# Row 15 (Murder): C/D/E become "0"/"0"/"***.*" text placeholders;
# F/G/L stay numeric with refreshed figures.
# ---------------------------------------------------------------------
Set-TextPlaceholder "C15" "0" "C22"
Set-TextPlaceholder "D15" "0" "D22"
Set-TextPlaceholder "E15" "***.*" "E22"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("L15").Value = -40

# ---------------------------------------------------------------------
# Row 16 (Rape): C16 flips from the "0" text placeholder back to a
# real number; the rest are refreshed figures.
# ---------------------------------------------------------------------
Set-NumberFromText "C16" 2 "D16"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 12.903225806451
$ws.Range("L16").Value = 34.615384615384

# ---------------------------------------------------------------------
# Row 17 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 127.272727272727
$ws.Range("I17").Value = 196
$ws.Range("J17").Value = 179
$ws.Range("K17").Value = 9.497206703910
$ws.Range("L17").Value = 53.125

# ---------------------------------------------------------------------
# Row 18 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = -24.657534246575
$ws.Range("L18").Value = -19.117647058823

# ---------------------------------------------------------------------
# Row 19 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 62.5
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 46.428571428571
$ws.Range("I19").Value = 348
$ws.Range("J19").Value = 335
$ws.Range("K19").Value = 3.880597014925
$ws.Range("L19").Value = 37.549407114624

# ---------------------------------------------------------------------
# Row 20 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = 8.641975308641
$ws.Range("L20").Value = 109.52380952381

# ---------------------------------------------------------------------
# Row 21 (G.L.A.) - bold TOTAL-style row, numbers only
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 47.058823529411
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = 25.714285714285
$ws.Range("I21").Value = 770
$ws.Range("J21").Value = 740
$ws.Range("K21").Value = 4.054054054054
$ws.Range("L21").Value = 37.254901960784

# ---------------------------------------------------------------------
# Row 23 (Transit): C/D/E become "0"/"0"/"***.*" text placeholders.
# ---------------------------------------------------------------------
Set-TextPlaceholder "C23" "0" "C22"
Set-TextPlaceholder "D23" "0" "D22"
Set-TextPlaceholder "E23" "***.*" "E22"

# ---------------------------------------------------------------------
# Row 24 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 86.363636363636
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 21.367521367521
$ws.Range("I24").Value = 1188
$ws.Range("J24").Value = 1141
$ws.Range("K24").Value = 4.119193689745
$ws.Range("L24").Value = 60.323886639676

# ---------------------------------------------------------------------
# Row 25 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -64.705882352941
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -26.190476190476
$ws.Range("I25").Value = 422
$ws.Range("J25").Value = 408
$ws.Range("K25").Value = 3.431372549019
$ws.Range("L25").Value = 30.650154798761

# ---------------------------------------------------------------------
# Row 26 (Misd. Assault): C/D/E become "0"/"0"/"***.*" text placeholders.
# ---------------------------------------------------------------------
Set-TextPlaceholder "C26" "0" "C22"
Set-TextPlaceholder "D26" "0" "D22"
Set-TextPlaceholder "E26" "***.*" "E22"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("L26").Value = -20.833333333333

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*): C/D/E become "0"/"0"/"***.*" text placeholders.
# ---------------------------------------------------------------------
Set-TextPlaceholder "C27" "0" "C22"
Set-TextPlaceholder "D27" "0" "D22"
Set-TextPlaceholder "E27" "***.*" "E22"
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 37.837837837837

# ---------------------------------------------------------------------
# Row 30 (Hate Crimes): G/H become "0"/"***.*" text placeholders
# (F30 is unchanged).
# ---------------------------------------------------------------------
Set-TextPlaceholder "G30" "0" "G22"
Set-TextPlaceholder "H30" "***.*" "H22"
